$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells: Wins / Losses / Ties in columns AC, AD, AE (29, 30, 31)
$ws.Cells.Item(1, 29).Value = "Wins"
$ws.Cells.Item(1, 30).Value = "Losses"
$ws.Cells.Item(1, 31).Value = "Ties"

# Copy the header style (bold/border/centered) from an existing header cell (AB1) to the new headers
$ws.Range("AB1").Copy()
$ws.Range("AC1:AE1").PasteSpecial(-4122)  # xlPasteFormats

# Fill in the season record (Wins=83, Losses=79, Ties=0) for every data row (2 through 46)
for ($r = 2; $r -le 46; $r++) {
    $ws.Cells.Item($r, 29).Value = 83
    $ws.Cells.Item($r, 30).Value = 79
    $ws.Cells.Item($r, 31).Value = 0
}
